$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date strings in column A use "DD/MM/YYYY"; they need to become
# "DD-MM-YYYY". Some of these strings (day <= 12) are ambiguous and would
# otherwise be auto-parsed by Excel into date serial numbers, so force the
# cell to Text format before writing, then restore the default style so the
# cell ends up identical (no explicit style) to how it started.
function Set-DateText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-DateText "A3"  "28-07-2022"
Set-DateText "A4"  "01-08-2022"
Set-DateText "A5"  "04-08-2022"
Set-DateText "A6"  "08-08-2022"
Set-DateText "A7"  "11-08-2022"
Set-DateText "A8"  "15-08-2022"
Set-DateText "A9"  "18-08-2022"
Set-DateText "A10" "22-08-2022"
Set-DateText "A11" "25-08-2022"
Set-DateText "A12" "29-08-2022"
Set-DateText "A13" "01-09-2022"
Set-DateText "A14" "05-09-2022"
Set-DateText "A15" "08-09-2022"
Set-DateText "A16" "12-09-2022"
Set-DateText "A17" "15-09-2022"
Set-DateText "A18" "19-09-2022"
Set-DateText "A19" "22-09-2022"
Set-DateText "A20" "26-09-2022"
Set-DateText "A21" "29-09-2022"

# Update the attendance counters that changed.
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0
